# Update cryptos list with refreshed prices / volume(1h) percentages.
# (Two coin-name pairs were also re-ordered by the upstream source: rows 32/33,
#  44/45 and 48/49 swap their Coin/Link contents.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that LOOKS like a number (e.g. "1.00", "0.993")
# into a cell while keeping it stored as TEXT (matching the source data,
# which is inline-string, not numeric). A direct $cell.Value = "1.00" would
# be auto-coerced by Excel into the number 1. Routing it through a text
# formula + Paste Special (values only) preserves the literal text "1.00"
# without touching the cell's style.
function Set-TextValue($address, $text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $ws.Range($address).PasteSpecial(-4163)
    $scratch.ClearContents()
}

# --- Plain text / already-unambiguous-as-text values -----------------------
$ws.Range("D2").Value = "60.822.68"
$ws.Range("E2").Value = "  +3.40%  "
$ws.Range("D3").Value = "2.695.21"
$ws.Range("E3").Value = "  +2.40%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("E5").Value = "  +2.09%  "
$ws.Range("E6").Value = "  +1.79%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  +2.43%  "
$ws.Range("D9").Value = "2.719.26"
$ws.Range("E9").Value = "  +2.33%  "
$ws.Range("E10").Value = "  +5.44%  "
$ws.Range("E11").Value = "  +1.63%  "
$ws.Range("E12").Value = "  +1.53%  "
$ws.Range("E13").Value = "  +2.75%  "
$ws.Range("D14").Value = "3.172.04"
$ws.Range("E14").Value = "  +2.34%  "
$ws.Range("D15").Value = "60.810.95"
$ws.Range("E15").Value = "  +3.45%  "
$ws.Range("D16").Value = "2.979.96"
$ws.Range("E16").Value = "  +12.44%  "
$ws.Range("E17").Value = "  +2.76%  "
$ws.Range("E18").Value = "  +1.46%  "
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("E21").Value = "  +2.98%  "
$ws.Range("E22").Value = "  +3.21%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("E24").Value = "  +3.32%  "
$ws.Range("E25").Value = "  +1.16%  "
$ws.Range("E26").Value = "  +5.86%  "
$ws.Range("E27").Value = "  -0.54%  "
$ws.Range("E28").Value = "  +4.43%  "
$ws.Range("E29").Value = "  +2.84%  "
$ws.Range("E30").Value = "  +10.38%  "
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("E32").Value = "  +2.16%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("E33").Value = "  +0.95%  "
$ws.Range("E34").Value = "  +0.35%  "
$ws.Range("E35").Value = "  +7.38%  "
$ws.Range("E36").Value = "  +10.34%  "
$ws.Range("E37").Value = "  -1.79%  "
$ws.Range("E38").Value = "  +5.74%  "
$ws.Range("E40").Value = "  +1.07%  "
$ws.Range("E41").Value = "  -0.57%  "
$ws.Range("E42").Value = "  +3.23%  "
$ws.Range("E43").Value = "  +3.47%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("E44").Value = "  +0.80%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("E45").Value = "  +0.74%  "
$ws.Range("D46").Value = "2.151.23"
$ws.Range("E46").Value = "  +8.59%  "
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E48").Value = "  +6.08%  "
$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("E49").Value = "  +2.44%  "
$ws.Range("E50").Value = "  +2.76%  "
$ws.Range("E51").Value = "  +1.74%  "

# --- Numeric-looking Price values that must remain TEXT ---------------------
Set-TextValue "D4" "1.00"
Set-TextValue "D5" "526.10"
Set-TextValue "D6" "145.88"
Set-TextValue "D10" "6.55"
Set-TextValue "D17" "21.42"
Set-TextValue "D19" "350.88"
Set-TextValue "D20" "4.54"
Set-TextValue "D21" "10.62"
Set-TextValue "D22" "6.35"
Set-TextValue "D24" "63.80"
Set-TextValue "D27" "0.993"
Set-TextValue "D32" "1.60"
Set-TextValue "D33" "19.13"
Set-TextValue "D34" "150.31"
Set-TextValue "D37" "0.949"
Set-TextValue "D40" "37.01"
Set-TextValue "D41" "3.68"
Set-TextValue "D42" "286.80"
Set-TextValue "D43" "20.20"
Set-TextValue "D44" "0.614"
Set-TextValue "D45" "0.0992"
Set-TextValue "D47" "0.996"
Set-TextValue "D48" "4.93"
Set-TextValue "D49" "0.0542"
